# Daily refresh of the cryptos list (prices/volumes) - GitHub Actions update.
# Numeric-looking price strings are forced to Text format before assignment so
# that Excel keeps them as literal strings (e.g. "1.00", "23.40") instead of
# silently coercing them to numbers and dropping trailing zeros.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.984.51"
$ws.Range("E2").Value = "  -0.59%  "
$ws.Range("D3").Value = "3.100.17"
$ws.Range("E3").Value = "  -2.01%  "
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "566.42"
$ws.Range("E5").Value = "  -0.82%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "160.30"
$ws.Range("E6").Value = "  -4.73%  "
$ws.Range("E7").Value = "  +0.07%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.570"
$ws.Range("E8").Value = "  -5.62%  "
$ws.Range("D9").Value = "3.111.06"
$ws.Range("E9").Value = "  -2.33%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.116"
$ws.Range("E10").Value = "  -2.53%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.57"
$ws.Range("E11").Value = "  -3.61%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.376"
$ws.Range("E12").Value = "  -3.19%  "
$ws.Range("D13").Value = "3.652.30"
$ws.Range("E13").Value = "  -1.63%  "
$ws.Range("E14").Value = "  -2.31%  "
$ws.Range("D15").Value = "64.117.63"
$ws.Range("E15").Value = "  -0.51%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "24.51"
$ws.Range("E16").Value = "  -3.25%  "
$ws.Range("D17").Value = "3.114.11"
$ws.Range("E17").Value = "  -1.33%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.0000154"
$ws.Range("E18").Value = "  -2.17%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "403.65"
$ws.Range("E19").Value = "  -3.20%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "5.18"
$ws.Range("E20").Value = "  -2.28%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.28"
$ws.Range("E21").Value = "  -4.54%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.93"
$ws.Range("E22").Value = "  -2.99%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.00"
$ws.Range("E23").Value = "  -0.01%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "67.35"
$ws.Range("E24").Value = "  -3.26%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.476"
$ws.Range("E25").Value = "  -4.38%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.192"
$ws.Range("E26").Value = "  -5.58%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0000101"
$ws.Range("E27").Value = "  -2.19%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.97"
$ws.Range("E28").Value = "  +1.84%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.00"
$ws.Range("E29").Value = "  +0.32%  "
$ws.Range("E30").Value = "  +0.13%  "
$ws.Range("E31").Value = "  -2.07%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "21.06"
$ws.Range("E32").Value = "  -3.17%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "164.07"
$ws.Range("E33").Value = "  +5.38%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.84"
$ws.Range("E34").Value = "  -4.22%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.17"
$ws.Range("E35").Value = "  -2.90%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.11"
$ws.Range("E36").Value = "  -0.49%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.34"
$ws.Range("E37").Value = "  -1.63%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.64"
$ws.Range("E38").Value = "  -3.37%  "
$ws.Range("D39").Value = "2.581.76"
$ws.Range("E39").Value = "  -4.43%  "
$ws.Range("B40").Value = "Filecoin"
$ws.Range("C40").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "4.07"
$ws.Range("E40").Value = "  -3.34%  "
$ws.Range("B41").Value = "EnergySwap"
$ws.Range("C41").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "23.40"
$ws.Range("E41").Value = "  -2.46%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "38.00"
$ws.Range("E42").Value = "  -2.76%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.686"
$ws.Range("E43").Value = "  -4.36%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0611"
$ws.Range("E44").Value = "  -0.78%  "
$ws.Range("B45").Value = "VeChain"
$ws.Range("C45").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0253"
$ws.Range("E45").Value = "  -3.92%  "
$ws.Range("B46").Value = "RenderToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "5.13"
$ws.Range("E46").Value = "  -5.93%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "20.76"
$ws.Range("E47").Value = "  -2.51%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "281.08"
$ws.Range("E48").Value = "  -3.02%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.999"
$ws.Range("E49").Value = "  -0.09%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0967"
$ws.Range("E50").Value = "  -2.22%  "
$ws.Range("E51").Value = "  +0.17%  "
